$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.616.91"
$ws.Range("E2").Value = "  +2.54%  "
$ws.Range("D3").Value = "1.860.23"
$ws.Range("E3").Value = "  +1.83%  "
$ws.Range("D4").Value = "'0.9996"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'245.00"
$ws.Range("E5").Value = "  +1.91%  "
$ws.Range("D6").Value = "'0.6977"
$ws.Range("E6").Value = "  +0.90%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").Value = "'0.07722"
$ws.Range("E8").Value = "  +1.33%  "
$ws.Range("D9").Value = "'0.3058"
$ws.Range("E9").Value = "  +0.66%  "
$ws.Range("D10").Value = "'23.68"
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("D11").Value = "'0.07751"
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("D12").Value = "'5.161"
$ws.Range("E12").Value = "  +2.40%  "
$ws.Range("D13").Value = "1.859.46"
$ws.Range("E13").Value = "  +1.88%  "
$ws.Range("D14").Value = "'92.36"
$ws.Range("E14").Value = "  +2.38%  "
$ws.Range("D15").Value = "'0.6919"
$ws.Range("E15").Value = "  +2.85%  "
$ws.Range("D16").Value = "'6.562"
$ws.Range("E16").Value = "  +3.13%  "
$ws.Range("D17").Value = "29.599.89"
$ws.Range("E17").Value = "  +2.53%  "
$ws.Range("D18").Value = "'0.000008335"
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").Value = "2.105.91"
$ws.Range("E19").Value = "  +1.70%  "
$ws.Range("D20").Value = "'241.73"
$ws.Range("E20").Value = "  -0.14%  "
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").Value = "'0.9999"
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("D23").Value = "'7.614"
$ws.Range("E23").Value = "  +2.82%  "
$ws.Range("E25").Value = "  +2.54%  "
$ws.Range("E26").Value = "  +1.96%  "
$ws.Range("D27").Value = "'159.51"
$ws.Range("E27").Value = "  -0.86%  "
$ws.Range("D28").Value = "'18.29"
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("D29").Value = "'1.536"
$ws.Range("E29").Value = "  +0.25%  "
$ws.Range("D30").Value = "'4.250"
$ws.Range("E30").Value = "  +1.35%  "
$ws.Range("D31").Value = "'4.182"
$ws.Range("E31").Value = "  +1.39%  "
$ws.Range("D32").Value = "'1.194"
$ws.Range("E32").Value = "  +0.21%  "
$ws.Range("D33").Value = "'0.05091"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "'0.7775"
$ws.Range("E34").Value = "  +3.80%  "
$ws.Range("D35").Value = "'1.897"
$ws.Range("E35").Value = "  +4.71%  "
$ws.Range("E36").Value = "  +1.48%  "
$ws.Range("D37").Value = "'2.685"
$ws.Range("E37").Value = "  +0.06%  "
$ws.Range("D38").Value = "1.324.63"
$ws.Range("E38").Value = "  +10.67%  "
$ws.Range("E39").Value = "  +2.08%  "
$ws.Range("E40").Value = "  +2.30%  "
$ws.Range("D41").Value = "'0.9596"
$ws.Range("E41").Value = "  +3.59%  "
$ws.Range("D42").Value = "'106.35"
$ws.Range("E42").Value = "  -1.70%  "
$ws.Range("D43").Value = "'5.821"
$ws.Range("E43").Value = "  +11.69%  "
$ws.Range("D44").Value = "'0.9998"
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").Value = "'0.00000000126"
$ws.Range("E45").Value = "  +3.94%  "
$ws.Range("D46").Value = "'9.766"
$ws.Range("E46").Value = "  +3.37%  "
$ws.Range("D47").Value = "2.005.39"
$ws.Range("E47").Value = "  +1.58%  "
$ws.Range("D48").Value = "'0.5212"
$ws.Range("E48").Value = "  +0.96%  "
$ws.Range("D49").Value = "'1.782"
$ws.Range("E49").Value = "  +3.39%  "
$ws.Range("D50").Value = "'64.45"
$ws.Range("E50").Value = "  +4.22%  "
$ws.Range("D51").Value = "'6.976"
$ws.Range("E51").Value = "  +1.60%  "

$ws.Range("D4,D5,D6,D8,D9,D10,D11,D12,D14,D15,D16,D18,D20,D22,D23,D27,D28,D29,D30,D31,D32,D33,D34,D35,D37,D41,D42,D43,D44,D45,D46,D48,D49,D50,D51").ClearFormats()
